$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 07:18:41"
$ws.Range("O2").Value = "1.7 °C"
$ws.Range("E3").Value = "2026-02-24 07:18:43"
$ws.Range("H3").Value = "'39%"
$ws.Range("E4").Value = "2026-02-24 07:18:46"
$ws.Range("H4").Value = "'89%"
$ws.Range("J4").Value = "1022.2 hPa"
$ws.Range("N4").Value = "4.2 °C 6:41 TU"
$ws.Range("O4").Value = "6.6 °C"
$ws.Range("E5").Value = "2026-02-24 07:18:49"
$ws.Range("M5").Value = "5.3 °C 6:59 TU"
$ws.Range("O5").Value = "4.0 °C"
$ws.Range("E6").Value = "2026-02-24 07:18:51"
$ws.Range("J6").Value = "1021.8 hPa"
$ws.Range("E7").Value = "2026-02-24 07:18:54"
$ws.Range("H7").Value = "'80%"
$ws.Range("J7").Value = "1021.9 hPa"
$ws.Range("K7").Value = "0.0 MJ/m2"
$ws.Range("N7").Value = "10.6 °C 6:34 TU"
$ws.Range("O7").Value = "11.9 °C"
$ws.Range("E8").Value = "2026-02-24 07:18:56"
$ws.Range("H8").Value = "'46%"
$ws.Range("J8").Value = "1021.5 hPa"
$ws.Range("L8").Value = "33.1 km/h - 305º 6:57 TU"
$ws.Range("N8").Value = "13.0 °C 6:36 TU"
$ws.Range("O8").Value = "14.7 °C"
$ws.Range("E9").Value = "2026-02-24 07:18:59"
$ws.Range("E10").Value = "2026-02-24 07:19:01"
$ws.Range("K10").Value = "0.0 MJ/m2"
$ws.Range("E11").Value = "2026-02-24 07:19:04"
$ws.Range("N11").Value = "1.4 °C 6:31 TU"
$ws.Range("O11").Value = "2.4 °C"
$ws.Range("E12").Value = "2026-02-24 07:19:07"
$ws.Range("N12").Value = "3.3 °C 6:39 TU"
$ws.Range("O12").Value = "5.5 °C"
$ws.Range("E13").Value = "2026-02-24 07:19:09"
$ws.Range("J13").Value = "1030.0 hPa"
$ws.Range("O13").Value = "-1.8 °C"
$ws.Range("E14").Value = "2026-02-24 07:19:12"
$ws.Range("O14").Value = "8.8 °C"
$ws.Range("E15").Value = "2026-02-24 07:19:15"
$ws.Range("O15").Value = "5.5 °C"
$ws.Range("E16").Value = "2026-02-24 07:19:17"
$ws.Range("E17").Value = "2026-02-24 07:19:19"
$ws.Range("H17").Value = "'33%"
$ws.Range("K17").Value = "0.0 MJ/m2"
$ws.Range("E18").Value = "2026-02-24 07:19:22"
$ws.Range("N18").Value = "1.1 °C 6:45 TU"
$ws.Range("O18").Value = "2.7 °C"
$ws.Range("E19").Value = "2026-02-24 07:19:25"
$ws.Range("H19").Value = "'68%"
$ws.Range("E20").Value = "2026-02-24 07:19:28"
$ws.Range("H20").Value = "'41%"
$ws.Range("K20").Value = "0.0 MJ/m2"
$ws.Range("N20").Value = "-1.1 °C 6:43 TU"
$ws.Range("O20").Value = "0.7 °C"
$ws.Range("E21").Value = "2026-02-24 07:19:30"
$ws.Range("H21").Value = "'82%"
$ws.Range("K21").Value = "0.0 MJ/m2"
$ws.Range("L21").Value = "6.5 km/h - 298º 6:43 TU"
$ws.Range("O21").Value = "3.1 °C"
$ws.Range("E22").Value = "2026-02-24 07:19:33"
$ws.Range("E23").Value = "2026-02-24 07:19:36"
$ws.Range("E24").Value = "2026-02-24 07:19:39"
$ws.Range("N24").Value = "0.6 °C 6:41 TU"
$ws.Range("O24").Value = "2.8 °C"
$ws.Range("E25").Value = "2026-02-24 07:19:41"
$ws.Range("H25").Value = "'34%"
$ws.Range("K25").Value = "0.0 MJ/m2"
$ws.Range("E26").Value = "2026-02-24 07:19:44"
$ws.Range("H26").Value = "'51%"
$ws.Range("J26").Value = "1022.4 hPa"
$ws.Range("L26").Value = "13.3 km/h - 354º 6:51 TU"
$ws.Range("M26").Value = "9.3 °C 6:47 TU"
$ws.Range("O26").Value = "7.7 °C"
$ws.Range("E27").Value = "2026-02-24 07:19:46"
$ws.Range("H27").Value = "'36%"
$ws.Range("O27").Value = "4.2 °C"
$ws.Range("E28").Value = "2026-02-24 07:19:49"
$ws.Range("H28").Value = "'96%"
$ws.Range("O28").Value = "3.5 °C"
$ws.Range("E29").Value = "2026-02-24 07:19:52"
$ws.Range("I29").Value = "0.1 mm"
$ws.Range("K29").Value = "0.0 MJ/m2"
$ws.Range("N29").Value = "3.0 °C 6:54 TU"
$ws.Range("O29").Value = "4.7 °C"
$ws.Range("E30").Value = "2026-02-24 07:19:54"
$ws.Range("N30").Value = "7.4 °C 6:30 TU"
$ws.Range("O30").Value = "8.8 °C"
$ws.Range("E31").Value = "2026-02-24 07:19:57"
$ws.Range("H31").Value = "'53%"
$ws.Range("E32").Value = "2026-02-24 07:19:59"
$ws.Range("E33").Value = "2026-02-24 07:20:02"
$ws.Range("H33").Value = "'71%"
$ws.Range("N33").Value = "0.2 °C 6:42 TU"
$ws.Range("O33").Value = "1.8 °C"
$ws.Range("E34").Value = "2026-02-24 07:20:05"
$ws.Range("H34").Value = "'54%"
$ws.Range("E35").Value = "2026-02-24 07:20:07"
$ws.Range("J35").Value = "1024.3 hPa"
$ws.Range("K35").Value = "0.0 MJ/m2"
$ws.Range("E36").Value = "2026-02-24 07:20:10"
$ws.Range("N36").Value = "5.3 °C 6:42 TU"
$ws.Range("O36").Value = "7.7 °C"
$ws.Range("E37").Value = "2026-02-24 07:20:13"
$ws.Range("H37").Value = "'97%"
$ws.Range("J37").Value = "1027.2 hPa"
$ws.Range("O37").Value = "0.8 °C"
$ws.Range("E38").Value = "2026-02-24 07:20:15"
$ws.Range("K38").Value = "0.0 MJ/m2"
$ws.Range("O38").Value = "6.3 °C"
$ws.Range("E39").Value = "2026-02-24 07:20:18"
$ws.Range("K39").Value = "0.0 MJ/m2"
$ws.Range("N39").Value = "2.3 °C 6:30 TU"
$ws.Range("O39").Value = "4.7 °C"
$ws.Range("E40").Value = "2026-02-24 07:20:20"
$ws.Range("H40").Value = "'96%"
$ws.Range("N40").Value = "-0.5 °C 6:41 TU"
$ws.Range("O40").Value = "0.9 °C"
$ws.Range("E41").Value = "2026-02-24 07:20:23"
$ws.Range("H41").Value = "'81%"
$ws.Range("N41").Value = "4.3 °C 6:39 TU"
$ws.Range("O41").Value = "6.7 °C"
$ws.Range("E42").Value = "2026-02-24 07:20:25"
$ws.Range("N42").Value = "4.6 °C 6:47 TU"
$ws.Range("O42").Value = "6.5 °C"
$ws.Range("E43").Value = "2026-02-24 07:20:27"
$ws.Range("K43").Value = "0.0 MJ/m2"
$ws.Range("O43").Value = "3.8 °C"
$ws.Range("E44").Value = "2026-02-24 07:20:29"
$ws.Range("H44").Value = "'49%"
$ws.Range("E45").Value = "2026-02-24 07:20:32"
$ws.Range("E46").Value = "2026-02-24 07:20:35"
$ws.Range("N46").Value = "0.4 °C 6:59 TU"
$ws.Range("O46").Value = "2.0 °C"
